$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: partial marks for writing code partially
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = "Partial marks for writing code partially"

# Row 30: for no output
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = "For no output"

# Row 37: for getting exceptions (deduction changed from 0 to -2.5)
$ws.Range("E37").Value = -2.5
$ws.Range("F37").Value = "For getting exceptions"

# Selection state to match author's final view position
$ws.Range("F37").Select()
